$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 08:20"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 188592
$ws.Range("C4").Value = 62
$ws.Range("E4").Value = 177286

# Row 16: Austria -> Austria
$ws.Range("B16").Value = 10192
$ws.Range("C16").Value = 12
$ws.Range("E16").Value = 8969

# Row 21: Israel -> Israel
$ws.Range("B21").Value = 5591
$ws.Range("C21").Value = 233
$ws.Range("E21").Value = 5346
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 21

# Row 22: Australia -> Australia
$ws.Range("B22").Value = 4862
$ws.Range("C22").Value = 99
$ws.Range("E22").Value = 4496

# Row 23: Noruega -> Noruega
$ws.Range("B23").Value = 4643
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 4591

# Row 37: Pakistan -> Pakistan
$ws.Range("B37").Value = 2042
$ws.Range("C37").Value = 104
$ws.Range("E37").Value = 1934

# Row 67: Barein -> Lituania
$ws.Range("A67").Value = "Lituania"
$ws.Range("B67").Value = 581
$ws.Range("C67").Value = 44
$ws.Range("D67").Value = 7
$ws.Range("E67").Value = 566
$ws.Range("F67").Value = 27
$ws.Range("H67").Value = 8

# Row 68: Lituania -> Barein
$ws.Range("A68").Value = "Barein"
$ws.Range("B68").Value = 567
$ws.Range("D68").Value = 295
$ws.Range("E68").Value = 268
$ws.Range("F68").Value = 2
$ws.Range("H68").Value = 4

# Row 70: Hungria -> Hungria
$ws.Range("E70").Value = 465
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 20

# Row 73: Bulgaria -> Bulgaria
$ws.Range("D73").Value = 20
$ws.Range("E73").Value = 384

# Row 77: Eslovaquia -> Kazajistan
$ws.Range("A77").Value = "Kazajistan"
$ws.Range("B77").Value = 369
$ws.Range("C77").Value = 26
$ws.Range("D77").Value = 24
$ws.Range("E77").Value = 342
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 3

# Row 78: Kazajistan -> Eslovaquia
$ws.Range("A78").Value = "Eslovaquia"
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 3
$ws.Range("E78").Value = 360
$ws.Range("F78").Value = 1
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 0

# Row 103: Ghana -> Bielorrusia
$ws.Range("A103").Value = "Bielorrusia"
$ws.Range("B103").Value = 163
$ws.Range("C103").Value = 11
$ws.Range("D103").Value = 53
$ws.Range("E103").Value = 109
$ws.Range("F103").Value = 2
$ws.Range("H103").Value = 1

# Row 104: Bielorrusia -> Ghana
$ws.Range("A104").Value = "Ghana"
$ws.Range("B104").Value = 161
$ws.Range("D104").Value = 31
$ws.Range("E104").Value = 125
$ws.Range("F104").Value = 1
$ws.Range("H104").Value = 5

# Row 107: Nigeria -> Venezuela
$ws.Range("A107").Value = "Venezuela"
$ws.Range("B107").Value = 143
$ws.Range("C107").Value = 8
$ws.Range("D107").Value = 41
$ws.Range("E107").Value = 99
$ws.Range("F107").Value = 6
$ws.Range("H107").Value = 3

# Row 108: Venezuela -> Nigeria
$ws.Range("A108").Value = "Nigeria"
$ws.Range("B108").Value = 139
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 9
$ws.Range("E108").Value = 128
$ws.Range("F108").Value = 0
$ws.Range("H108").Value = 2

# Row 115: Montenegro -> Montenegro
$ws.Range("F115").Value = 4

# Row 116: Camboya -> Camboya
$ws.Range("D116").Value = 25
$ws.Range("E116").Value = 84

# Row 151: Maldivas -> Gabon
$ws.Range("A151").Value = "Gabon"
$ws.Range("C151").Value = 2
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 17
$ws.Range("H151").Value = 1

# Row 152: Islas Virgenes de los Estados Unidos -> Maldivas
$ws.Range("A152").Value = "Maldivas"
$ws.Range("B152").Value = 18
$ws.Range("D152").Value = 13
$ws.Range("E152").Value = 5

# Row 153: Nueva Caledonia -> Islas Virgenes de los Estados Unidos
$ws.Range("A153").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B153").Value = 17
$ws.Range("E153").Value = 17

# Row 154: Haiti -> Nueva Caledonia
$ws.Range("A154").Value = "Nueva Caledonia"
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 16

# Row 155: Gabon -> Haiti
$ws.Range("A155").Value = "Haiti"
$ws.Range("C155").Value = 1
$ws.Range("D155").Value = 1
$ws.Range("H155").Value = 0

# Row 162: Santa Lucia -> Mongolia
$ws.Range("A162").Value = "Mongolia"
$ws.Range("B162").Value = 14
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 2

# Row 163: Dominica -> Santa Lucia
$ws.Range("A163").Value = "Santa Lucia"
$ws.Range("B163").Value = 13
$ws.Range("D163").Value = 1

# Row 164: Mongolia -> Dominica
$ws.Range("A164").Value = "Dominica"
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 12
